# PSI_inputs.xlsx - "Minor fixes as running code"
#
# Adds a new soil-soil parameter row ("pipelay_Fct") to the Inputs sheet,
# right above the existing "int_SHANSEP_S" row (i.e. becomes the new row 46,
# pushing the previous rows 46-49 down to 47-50). Also nudges the saved
# worksheet view/selection on the Inputs tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")
$ws.Activate()

# --- Insert the new row, shifting int_SHANSEP_S / int_SHANSEP_m / delta / blank row down ---
$ws.Rows("46:46").Insert()

# --- Populate the new row 46 ---
$ws.Range("A46").Value = "pipelay_Fct"

$ws.Range("B46:D46").NumberFormat = "0"
$ws.Range("B46").Value = 10
$ws.Range("C46").Value = 50
$ws.Range("D46").Value = 100

$ws.Range("E46").NumberFormat = "0.0"
$ws.Range("E46").Value = 1

$ws.Range("F46").Value = "Automated Fit"

$ws.Range("G46").Value = "% soil-soil: factor applied to phi to account for dynamic effects during pipelay, similar to pipelay_St for undrained"

# --- Restore the view/selection state on the Inputs tab ---
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("P30").Select()

Write-Output "Inserted pipelay_Fct row and updated Inputs view"
